$wb = $excel.ActiveWorkbook

# ----- Sheet "VENTAS POR GRUPO" -----
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Widen column I (9th column) from 9 -> 12 characters.
# ColumnWidth is expressed in characters; 11.14 lands on the stored width of 12
# given this workbook's default font metrics.
$wsVentasGrupo.Columns.Item(9).ColumnWidth = 11.14

$wsVentasGrupo.Range("D4").Value = 890.88
$wsVentasGrupo.Range("L4").Value = 2787.16
$wsVentasGrupo.Range("L5").Value = 1249.97
$wsVentasGrupo.Range("I50").Value = 176.16
$wsVentasGrupo.Range("D53").Value = "4 de 51"
$wsVentasGrupo.Range("I53").Value = "1 de 51"

# ----- Sheet "VENTA MENSUAL" -----
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsVentaMensual.Range("F4").Value = 4175.7
$wsVentaMensual.Range("F5").Value = 1252.09
$wsVentaMensual.Range("F50").Value = 176.16
$wsVentaMensual.Range("F53").Value = 16304.17

# ----- Sheet "CUMPLIMIENTO MENSUAL" -----
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 3717.32
$wsCumplimiento.Range("E3").Value = 23739.6876
$wsCumplimiento.Range("F3").Value = 0.1353869312400962

$wsCumplimiento.Range("D9").Value = 176.16
$wsCumplimiento.Range("E9").Value = 123.84
$wsCumplimiento.Range("F9").Value = 0.5871999999999999

$wsCumplimiento.Range("D10").Value = 131.37
$wsCumplimiento.Range("E10").Value = 1169.13
$wsCumplimiento.Range("F10").Value = 0.1010149942329873

$wsCumplimiento.Range("D16").Value = 5225.81
$wsCumplimiento.Range("E16").Value = 27515.64
$wsCumplimiento.Range("F16").Value = 0.1596083863115409

$wsCumplimiento.Range("D19").Value = 16304.17
$wsCumplimiento.Range("E19").Value = 78143.27064517915
$wsCumplimiento.Range("F19").Value = 0.1726269117365671
